$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 0.3074763333333333
$ws.Range("H2").Value = 0.9224289999999999
$ws.Range("I2").Value = 0.1203396906281759
$ws.Range("J2").Value = 0.1203396906281759
$ws.Range("M2").Value = 1.646588666666666
$ws.Range("N2").Value = 4.939766
$ws.Range("O2").Value = 0.039310317935267
$ws.Range("P2").Value = 0.039310317935267
$ws.Range("Q2").Value = 0.5062870457348888
$ws.Range("R2").Value = 4.556583411614
$ws.Range("S2").Value = 0.004730591498825264
$ws.Range("T2").Value = 0.004730591498825264

# Row 3
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 0.3074763333333333
$ws.Range("H3").Value = 0.9224289999999999
$ws.Range("I3").Value = 0.1203396906281759
$ws.Range("J3").Value = 0.1203396906281759
$ws.Range("M3").Value = 11.67754066666667
$ws.Range("N3").Value = 35.032622
$ws.Range("O3").Value = 0.278787195370394
$ws.Range("P3").Value = 0.278787195370394
$ws.Range("Q3").Value = 3.590567386537555
$ws.Range("R3").Value = 32.315106478838
$ws.Range("S3").Value = 0.03354916484197003
$ws.Range("T3").Value = 0.03354916484197003

# Row 4
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 0.3074763333333333
$ws.Range("H4").Value = 0.9224289999999999
$ws.Range("I4").Value = 0.1203396906281759
$ws.Range("J4").Value = 0.1203396906281759
$ws.Range("M4").Value = 0.7553226666666667
$ws.Range("N4").Value = 2.265968
$ws.Range("O4").Value = 0.01803241742850595
$ws.Range("P4").Value = 0.01803241742850595
$ws.Range("Q4").Value = 0.2322438440302222
$ws.Range("R4").Value = 2.090194596272
$ws.Range("S4").Value = 0.002170015534624533
$ws.Range("T4").Value = 0.002170015534624533

# Row 5
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.3074763333333333
$ws.Range("H5").Value = 0.9224289999999999
$ws.Range("I5").Value = 0.1203396906281759
$ws.Range("J5").Value = 0.1203396906281759
$ws.Range("M5").Value = 27.21325766666666
$ws.Range("N5").Value = 81.63977299999999
$ws.Range("O5").Value = 0.6496836961088899
$ws.Range("P5").Value = 0.6496836961088899
$ws.Range("Q5").Value = 8.367432685401887
$ws.Range("R5").Value = 75.30689416861699
$ws.Range("S5").Value = 0.07818273499591363
$ws.Range("T5").Value = 0.07818273499591363

# Row 6
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.3074763333333333
$ws.Range("H6").Value = 0.9224289999999999
$ws.Range("I6").Value = 0.1203396906281759
$ws.Range("J6").Value = 0.1203396906281759
$ws.Range("M6").Value = 0.5942236666666667
$ws.Range("N6").Value = 1.782671
$ws.Range("O6").Value = 0.01418637315694314
$ws.Range("P6").Value = 0.01418637315694314
$ws.Range("Q6").Value = 0.1827097142065555
$ws.Range("R6").Value = 1.644387427859
$ws.Range("S6").Value = 0.001707183756842396
$ws.Range("T6").Value = 0.001707183756842396

# Row 7
$ws.Range("I7").Value = 0.1296640274695671
$ws.Range("J7").Value = 0.129664027469567
$ws.Range("M7").Value = 1.646588666666666
$ws.Range("N7").Value = 4.939766
$ws.Range("O7").Value = 0.039310317935267
$ws.Range("P7").Value = 0.039310317935267
$ws.Range("Q7").Value = 0.5455159229924444
$ws.Range("R7").Value = 4.909643306932
$ws.Range("S7").Value = 0.005097134144595875
$ws.Range("T7").Value = 0.005097134144595874

# Row 8
$ws.Range("I8").Value = 0.1296640274695671
$ws.Range("J8").Value = 0.129664027469567
$ws.Range("O8").Value = 0.278787195370394
$ws.Range("P8").Value = 0.278787195370394
$ws.Range("S8").Value = 0.03614867055867033
$ws.Range("T8").Value = 0.03614867055867032

# Row 9
$ws.Range("I9").Value = 0.1296640274695671
$ws.Range("J9").Value = 0.129664027469567
$ws.Range("M9").Value = 0.7553226666666667
$ws.Range("N9").Value = 2.265968
$ws.Range("O9").Value = 0.01803241742850595
$ws.Range("P9").Value = 0.01803241742850595
$ws.Range("Q9").Value = 0.2502389030151111
$ws.Range("R9").Value = 2.252150127136
$ws.Range("S9").Value = 0.002338155868792495
$ws.Range("T9").Value = 0.002338155868792495

# Row 10
$ws.Range("I10").Value = 0.1296640274695671
$ws.Range("J10").Value = 0.129664027469567
$ws.Range("M10").Value = 27.21325766666666
$ws.Range("N10").Value = 81.63977299999999
$ws.Range("O10").Value = 0.6496836961088899
$ws.Range("P10").Value = 0.6496836961088899
$ws.Range("Q10").Value = 9.015770407138444
$ws.Range("R10").Value = 81.141933664246
$ws.Range("S10").Value = 0.08424060461879296
$ws.Range("T10").Value = 0.08424060461879294

# Row 11
$ws.Range("I11").Value = 0.1296640274695671
$ws.Range("J11").Value = 0.129664027469567
$ws.Range("M11").Value = 0.5942236666666667
$ws.Range("N11").Value = 1.782671
$ws.Range("O11").Value = 0.01418637315694314
$ws.Range("P11").Value = 0.01418637315694314
$ws.Range("Q11").Value = 0.1968666969157778
$ws.Range("R11").Value = 1.771800272242
$ws.Range("S11").Value = 0.001839462278715404
$ws.Range("T11").Value = 0.001839462278715403

# Row 12
$ws.Range("G12").Value = 0.4806033333333333
$ws.Range("H12").Value = 1.44181
$ws.Range("I12").Value = 0.1880979125164216
$ws.Range("J12").Value = 0.1880979125164216
$ws.Range("M12").Value = 1.646588666666666
$ws.Range("N12").Value = 4.939766
$ws.Range("O12").Value = 0.039310317935267
$ws.Range("P12").Value = 0.039310317935267
$ws.Range("Q12").Value = 0.7913560018288887
$ws.Range("R12").Value = 7.122204016459999
$ws.Range("S12").Value = 0.007394188743980571
$ws.Range("T12").Value = 0.00739418874398057

# Row 13
$ws.Range("G13").Value = 0.4806033333333333
$ws.Range("H13").Value = 1.44181
$ws.Range("I13").Value = 0.1880979125164216
$ws.Range("J13").Value = 0.1880979125164216
$ws.Range("O13").Value = 0.278787195370394
$ws.Range("P13").Value = 0.278787195370394
$ws.Range("Q13").Value = 5.612264969535556
$ws.Range("R13").Value = 50.51038472582
$ws.Range("S13").Value = 0.05243928948547891
$ws.Range("T13").Value = 0.0524392894854789

# Row 14
$ws.Range("G14").Value = 0.4806033333333333
$ws.Range("H14").Value = 1.44181
$ws.Range("I14").Value = 0.1880979125164216
$ws.Range("J14").Value = 0.1880979125164216
$ws.Range("M14").Value = 0.7553226666666667
$ws.Range("N14").Value = 2.265968
$ws.Range("O14").Value = 0.01803241742850595
$ws.Range("P14").Value = 0.01803241742850595
$ws.Range("Q14").Value = 0.3630105913422222
$ws.Range("R14").Value = 3.267095322079999
$ws.Range("S14").Value = 0.003391860075926708
$ws.Range("T14").Value = 0.003391860075926708

# Row 15
$ws.Range("G15").Value = 0.4806033333333333
$ws.Range("H15").Value = 1.44181
$ws.Range("I15").Value = 0.1880979125164216
$ws.Range("J15").Value = 0.1880979125164216
$ws.Range("M15").Value = 27.21325766666666
$ws.Range("N15").Value = 81.63977299999999
$ws.Range("O15").Value = 0.6496836961088899
$ws.Range("P15").Value = 0.6496836961088899
$ws.Range("Q15").Value = 13.07878234545888
$ws.Range("R15").Value = 117.70904110913
$ws.Range("S15").Value = 0.1222041470340354
$ws.Range("T15").Value = 0.1222041470340354

# Row 16
$ws.Range("G16").Value = 0.4806033333333333
$ws.Range("H16").Value = 1.44181
$ws.Range("I16").Value = 0.1880979125164216
$ws.Range("J16").Value = 0.1880979125164216
$ws.Range("M16").Value = 0.5942236666666667
$ws.Range("N16").Value = 1.782671
$ws.Range("O16").Value = 0.01418637315694314
$ws.Range("P16").Value = 0.01418637315694314
$ws.Range("Q16").Value = 0.2855858749455555
$ws.Range("R16").Value = 2.57027287451
$ws.Range("S16").Value = 0.002668427177000002
$ws.Range("T16").Value = 0.002668427177000002

# Row 17
$ws.Range("G17").Value = 0.225236
$ws.Range("H17").Value = 0.6757080000000001
$ws.Range("I17").Value = 0.08815257507622103
$ws.Range("J17").Value = 0.08815257507622101
$ws.Range("M17").Value = 1.646588666666666
$ws.Range("N17").Value = 4.939766
$ws.Range("O17").Value = 0.039310317935267
$ws.Range("P17").Value = 0.039310317935267
$ws.Range("Q17").Value = 0.3708710449253333
$ws.Range("R17").Value = 3.337839404328
$ws.Range("S17").Value = 0.003465305753058742
$ws.Range("T17").Value = 0.003465305753058742

# Row 18
$ws.Range("G18").Value = 0.225236
$ws.Range("H18").Value = 0.6757080000000001
$ws.Range("I18").Value = 0.08815257507622103
$ws.Range("J18").Value = 0.08815257507622101
$ws.Range("O18").Value = 0.278787195370394
$ws.Range("P18").Value = 0.278787195370394
$ws.Range("Q18").Value = 2.630202549597334
$ws.Range("R18").Value = 23.671822946376
$ws.Range("S18").Value = 0.02457580917017775
$ws.Range("T18").Value = 0.02457580917017775

# Row 19
$ws.Range("G19").Value = 0.225236
$ws.Range("H19").Value = 0.6757080000000001
$ws.Range("I19").Value = 0.08815257507622103
$ws.Range("J19").Value = 0.08815257507622101
$ws.Range("M19").Value = 0.7553226666666667
$ws.Range("N19").Value = 2.265968
$ws.Range("O19").Value = 0.01803241742850595
$ws.Range("P19").Value = 0.01803241742850595
$ws.Range("Q19").Value = 0.1701258561493333
$ws.Range("R19").Value = 1.531132705344
$ws.Range("S19").Value = 0.001589604031172127
$ws.Range("T19").Value = 0.001589604031172127

# Row 20
$ws.Range("G20").Value = 0.225236
$ws.Range("H20").Value = 0.6757080000000001
$ws.Range("I20").Value = 0.08815257507622103
$ws.Range("J20").Value = 0.08815257507622101
$ws.Range("M20").Value = 27.21325766666666
$ws.Range("N20").Value = 81.63977299999999
$ws.Range("O20").Value = 0.6496836961088899
$ws.Range("P20").Value = 0.6496836961088899
$ws.Range("Q20").Value = 6.129405303809333
$ws.Range("R20").Value = 55.164647734284
$ws.Range("S20").Value = 0.05727129079703568
$ws.Range("T20").Value = 0.05727129079703567

# Row 21
$ws.Range("G21").Value = 0.225236
$ws.Range("H21").Value = 0.6757080000000001
$ws.Range("I21").Value = 0.08815257507622103
$ws.Range("J21").Value = 0.08815257507622101
$ws.Range("M21").Value = 0.5942236666666667
$ws.Range("N21").Value = 1.782671
$ws.Range("O21").Value = 0.01418637315694314
$ws.Range("P21").Value = 0.01418637315694314
$ws.Range("Q21").Value = 0.1338405617853334
$ws.Range("R21").Value = 1.204565056068
$ws.Range("S21").Value = 0.001250565324776717
$ws.Range("T21").Value = 0.001250565324776717

# Row 22
$ws.Range("G22").Value = 1.210453666666667
$ws.Range("H22").Value = 3.631361
$ws.Range("I22").Value = 0.4737457943096146
$ws.Range("J22").Value = 0.4737457943096145
$ws.Range("M22").Value = 1.646588666666666
$ws.Range("N22").Value = 4.939766
$ws.Range("O22").Value = 0.039310317935267
$ws.Range("P22").Value = 0.039310317935267
$ws.Range("Q22").Value = 1.993119289058444
$ws.Range("R22").Value = 17.938073601526
$ws.Range("S22").Value = 0.01862309779480655
$ws.Range("T22").Value = 0.01862309779480655

# Row 23
$ws.Range("G23").Value = 1.210453666666667
$ws.Range("H23").Value = 3.631361
$ws.Range("I23").Value = 0.4737457943096146
$ws.Range("J23").Value = 0.4737457943096145
$ws.Range("O23").Value = 0.278787195370394
$ws.Range("P23").Value = 0.278787195370394
$ws.Range("Q23").Value = 14.13512191761578
$ws.Range("R23").Value = 127.216097258542
$ws.Range("S23").Value = 0.132074261314097
$ws.Range("T23").Value = 0.132074261314097

# Row 24
$ws.Range("G24").Value = 1.210453666666667
$ws.Range("H24").Value = 3.631361
$ws.Range("I24").Value = 0.4737457943096146
$ws.Range("J24").Value = 0.4737457943096145
$ws.Range("M24").Value = 0.7553226666666667
$ws.Range("N24").Value = 2.265968
$ws.Range("O24").Value = 0.01803241742850595
$ws.Range("P24").Value = 0.01803241742850595
$ws.Range("Q24").Value = 0.9142830913831113
$ws.Range("R24").Value = 8.228547822448
$ws.Range("S24").Value = 0.008542781917990088
$ws.Range("T24").Value = 0.008542781917990086

# Row 25
$ws.Range("G25").Value = 1.210453666666667
$ws.Range("H25").Value = 3.631361
$ws.Range("I25").Value = 0.4737457943096146
$ws.Range("J25").Value = 0.4737457943096145
$ws.Range("M25").Value = 27.21325766666666
$ws.Range("N25").Value = 81.63977299999999
$ws.Range("O25").Value = 0.6496836961088899
$ws.Range("P25").Value = 0.6496836961088899
$ws.Range("Q25").Value = 32.94038752456144
$ws.Range("R25").Value = 296.463487721053
$ws.Range("S25").Value = 0.3077849186631123
$ws.Range("T25").Value = 0.3077849186631123

# Row 26
$ws.Range("G26").Value = 1.210453666666667
$ws.Range("H26").Value = 3.631361
$ws.Range("I26").Value = 0.4737457943096146
$ws.Range("J26").Value = 0.4737457943096145
$ws.Range("M26").Value = 0.5942236666666667
$ws.Range("N26").Value = 1.782671
$ws.Range("O26").Value = 0.01418637315694314
$ws.Range("P26").Value = 0.01418637315694314
$ws.Range("Q26").Value = 0.7192802161367778
$ws.Range("R26").Value = 6.473521945231001
$ws.Range("S26").Value = 0.006720734619608621
$ws.Range("T26").Value = 0.00672073461960862
